$d = $word.ActiveDocument

# 1. Append the new "next steps" section after the very last paragraph
#    (the one holding the final screenshot), right before the sectPr.

$origCount = $d.Paragraphs.Count

# Find an existing bold "Báo cáo tuần N:" heading paragraph so we can copy
# its full formatting (run rPr AND paragraph-mark rPr, i.e. b + bCs on both)
# onto our brand new heading paragraph.
$headingSrc = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*cáo tuần 3*") {
        $headingSrc = $cand
        break
    }
}

# Insert a new paragraph at the end of the document, then stamp it with the
# fully-formatted (bold, bCs) text copied from the source heading. Copying
# the *entire* paragraph range (including its own trailing paragraph mark)
# is what makes Word replicate the bold onto the new paragraph mark too;
# it also leaves one extra empty paragraph behind, which we reuse below.
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Content
$r.Collapse(0)
$r.FormattedText = $headingSrc.Range.FormattedText

# Rewrite the copied paragraph's text (keeping the bold formatting intact).
$headingPara = $d.Paragraphs.Item($origCount + 1)
$headingTextRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$headingTextRange.Text = "Công việc thực hiện tiếp theo vào tuần 5:"

# The FormattedText assignment above left a trailing empty (non-bold)
# paragraph at the end of the document — use it for the first bullet line.
$r = $d.Content
$r.Collapse(0)
$r.Text = "- Thực hiện chức năng CRUD cho các mục có trong website."

# Remaining bullet paragraphs.
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Content
$r.Collapse(0)
$r.Text = "- Liên kết dữ liệu giữa fontend và backend."

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()
$r = $d.Content
$r.Collapse(0)
$r.Text = "- Tạo database."

# 2. Mark every inline picture's run as NoProofing (w:noProof) — matches
#    the 6 <w:noProof/> additions next to the image runs in the diff. Done
#    last so it cannot bleed into the freshly typed paragraphs above.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = 1
}
